# Ran Priors with _1
# Update the prior "variance-related" column D and column E for rows 2-44:
#   D: 1 -> 5
#   E: 0.15 or 0.05 -> 0.2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 4).Value = 5
    $ws.Cells.Item($r, 5).Value = 0.2
}

# Update the active view/selection to match the author's final state
$ws.Range("E2:E44").Select()
$ws.Application.ActiveWindow.ScrollRow = 25
